$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(53, 8).Value = 317.44446
$ws.Cells.Item(53, 9).Value = 355.7143
$ws.Cells.Item(53, 11).Value = 355.7143
$ws.Cells.Item(53, 13).Value = 281.2857

$ws.Cells.Item(132, 8).Value = 4001.9333
$ws.Cells.Item(132, 9).Value = 4129.0713
$ws.Cells.Item(132, 11).Value = 12387.2139
$ws.Cells.Item(132, 13).Value = -9857.213899999999

$ws.Cells.Item(135, 8).Value = 614.25
$ws.Cells.Item(135, 9).Value = 625.5263
$ws.Cells.Item(135, 11).Value = 5629.736699999999
$ws.Cells.Item(135, 13).Value = -3094.736699999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 2064.5454
$ws.Cells.Item(74, 9).Value = 2071.6
$ws.Cells.Item(74, 11).Value = 2071.6
$ws.Cells.Item(74, 13).Value = -1197.6

$ws.Cells.Item(77, 8).Value = 2064.5454
$ws.Cells.Item(77, 9).Value = 2071.6
$ws.Cells.Item(77, 11).Value = 10358
$ws.Cells.Item(77, 13).Value = -5990

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(35, 8).Value = 65036.5
$ws.Cells.Item(35, 10).Value = 65036.5
$ws.Cells.Item(35, 12).Value = 65036.5
$ws.Cells.Item(35, 14).Value = -65656.5

$ws.Cells.Item(86, 8).Value = 7577
$ws.Cells.Item(86, 10).Value = 19129.6
$ws.Cells.Item(86, 12).Value = 19129.6
$ws.Cells.Item(86, 14).Value = -21375.6

$ws.Cells.Item(89, 8).Value = 7577
$ws.Cells.Item(89, 10).Value = 19129.6
$ws.Cells.Item(89, 12).Value = 95648
$ws.Cells.Item(89, 14).Value = -106880

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(20, 8).Value = 68326.336
$ws.Cells.Item(20, 10).Value = 68326.336
$ws.Cells.Item(20, 12).Value = 68326.336
$ws.Cells.Item(20, 14).Value = -68798.336

$ws.Cells.Item(30, 8).Value = 68326.336
$ws.Cells.Item(30, 10).Value = 68326.336
$ws.Cells.Item(30, 12).Value = 68326.336
$ws.Cells.Item(30, 14).Value = -68508.336

$ws.Cells.Item(119, 8).Value = 40000
$ws.Cells.Item(119, 10).Value = 40000
$ws.Cells.Item(119, 12).Value = 40000
$ws.Cells.Item(119, 14).Value = -49676

$ws.Cells.Item(128, 8).Value = 68326.336
$ws.Cells.Item(128, 10).Value = 68326.336
$ws.Cells.Item(128, 12).Value = 68326.336
$ws.Cells.Item(128, 14).Value = -78286.336

$ws.Cells.Item(132, 8).Value = 2645.625
$ws.Cells.Item(132, 9).Value = 2595.5715
$ws.Cells.Item(132, 11).Value = 7786.7145
$ws.Cells.Item(132, 13).Value = -5256.7145

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 2660.077
$ws.Cells.Item(5, 9).Value = 2757
$ws.Cells.Item(5, 10).Value = 1497
$ws.Cells.Item(5, 11).Value = 8271
$ws.Cells.Item(5, 12).Value = 4491
$ws.Cells.Item(5, 13).Value = -8159
$ws.Cells.Item(5, 14).Value = -4715

$ws.Cells.Item(113, 8).Value = 4071.2273
$ws.Cells.Item(113, 9).Value = 3865.8333
$ws.Cells.Item(113, 10).Value = 4317.7
$ws.Cells.Item(113, 11).Value = 11597.4999
$ws.Cells.Item(113, 12).Value = 12953.1
$ws.Cells.Item(113, 13).Value = -9427.499899999999
$ws.Cells.Item(113, 14).Value = -17293.1

$ws.Cells.Item(135, 8).Value = 2660.077
$ws.Cells.Item(135, 9).Value = 2757
$ws.Cells.Item(135, 10).Value = 1497
$ws.Cells.Item(135, 11).Value = 24813
$ws.Cells.Item(135, 12).Value = 13473
$ws.Cells.Item(135, 13).Value = -22278
$ws.Cells.Item(135, 14).Value = -18543

$ws.Cells.Item(137, 8).Value = 1997
$ws.Cells.Item(137, 9).Value = 1999
$ws.Cells.Item(137, 11).Value = 5997
$ws.Cells.Item(137, 13).Value = -897

$ws.Cells.Item(139, 8).Value = 6888.909
$ws.Cells.Item(139, 9).Value = 7077.8
$ws.Cells.Item(139, 10).Value = 5000
$ws.Cells.Item(139, 11).Value = 21233.4
$ws.Cells.Item(139, 12).Value = 15000
$ws.Cells.Item(139, 13).Value = -16093.4
$ws.Cells.Item(139, 14).Value = -25280

$ws.Cells.Item(140, 8).Value = 670921.7
$ws.Cells.Item(140, 9).Value = 670921.7
$ws.Cells.Item(140, 11).Value = 2012765.1
$ws.Cells.Item(140, 13).Value = -2007585.1

$ws.Cells.Item(141, 8).Value = 8515
$ws.Cells.Item(141, 9).Value = 8515
$ws.Cells.Item(141, 11).Value = 25545
$ws.Cells.Item(141, 13).Value = -20365

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 2841.9
$ws.Cells.Item(80, 9).Value = 3834.6667
$ws.Cells.Item(80, 10).Value = 2416.4285
$ws.Cells.Item(80, 11).Value = 3834.6667
$ws.Cells.Item(80, 12).Value = 2416.4285
$ws.Cells.Item(80, 13).Value = -2836.6667
$ws.Cells.Item(80, 14).Value = -4412.4285

$ws.Cells.Item(83, 8).Value = 2841.9
$ws.Cells.Item(83, 9).Value = 3834.6667
$ws.Cells.Item(83, 10).Value = 2416.4285
$ws.Cells.Item(83, 11).Value = 19173.3335
$ws.Cells.Item(83, 12).Value = 12082.1425
$ws.Cells.Item(83, 13).Value = -14181.3335
$ws.Cells.Item(83, 14).Value = -22066.1425

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 745.6
$ws.Cells.Item(22, 9).Value = 682
$ws.Cells.Item(22, 10).Value = 1000
$ws.Cells.Item(22, 11).Value = 682
$ws.Cells.Item(22, 12).Value = 1000
$ws.Cells.Item(22, 13).Value = -387
$ws.Cells.Item(22, 14).Value = -1590

$ws.Cells.Item(27, 8).Value = 745.6
$ws.Cells.Item(27, 9).Value = 682
$ws.Cells.Item(27, 10).Value = 1000
$ws.Cells.Item(27, 11).Value = 682
$ws.Cells.Item(27, 12).Value = 1000
$ws.Cells.Item(27, 13).Value = -575
$ws.Cells.Item(27, 14).Value = -1214

$ws.Cells.Item(55, 8).Value = 374.1111
$ws.Cells.Item(55, 9).Value = 236.33333
$ws.Cells.Item(55, 10).Value = 649.6667
$ws.Cells.Item(55, 11).Value = 236.33333
$ws.Cells.Item(55, 12).Value = 649.6667
$ws.Cells.Item(55, 13).Value = -63.33332999999999
$ws.Cells.Item(55, 14).Value = -995.6667

$ws.Cells.Item(134, 8).Value = 84999.10000000001
$ws.Cells.Item(134, 10).Value = 84999.10000000001
$ws.Cells.Item(134, 12).Value = 84999.10000000001
$ws.Cells.Item(134, 14).Value = -95139.10000000001

$ws.Cells.Item(136, 8).Value = 7817.6875
$ws.Cells.Item(136, 9).Value = 6237.154
$ws.Cells.Item(136, 11).Value = 18711.462
$ws.Cells.Item(136, 13).Value = -16161.462

$ws.Cells.Item(139, 8).Value = 75000
$ws.Cells.Item(139, 10).Value = 75000
$ws.Cells.Item(139, 12).Value = 75000
$ws.Cells.Item(139, 14).Value = -85280

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 5349.2
$ws.Cells.Item(81, 9).Value = 5588
$ws.Cells.Item(81, 10).Value = 3200
$ws.Cells.Item(81, 11).Value = 11176
$ws.Cells.Item(81, 12).Value = 6400
$ws.Cells.Item(81, 13).Value = -10115
$ws.Cells.Item(81, 14).Value = -8522

$ws.Cells.Item(84, 8).Value = 5349.2
$ws.Cells.Item(84, 9).Value = 5588
$ws.Cells.Item(84, 10).Value = 3200
$ws.Cells.Item(84, 11).Value = 55880
$ws.Cells.Item(84, 12).Value = 32000
$ws.Cells.Item(84, 13).Value = -50576
$ws.Cells.Item(84, 14).Value = -42608

$ws.Cells.Item(132, 8).Value = 2576
$ws.Cells.Item(132, 9).Value = 2363.8857
$ws.Cells.Item(132, 10).Value = 10000
$ws.Cells.Item(132, 11).Value = 7091.657099999999
$ws.Cells.Item(132, 12).Value = 30000
$ws.Cells.Item(132, 13).Value = -4561.657099999999
$ws.Cells.Item(132, 14).Value = -35060

$ws.Cells.Item(136, 8).Value = 11394.467
$ws.Cells.Item(136, 9).Value = 9642.583000000001
$ws.Cells.Item(136, 10).Value = 18402
$ws.Cells.Item(136, 11).Value = 28927.749
$ws.Cells.Item(136, 12).Value = 55206
$ws.Cells.Item(136, 13).Value = -26377.749
$ws.Cells.Item(136, 14).Value = -60306
